$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps text formatting, matching the source data
# (values like "0.7525" would otherwise be auto-converted to numbers by Excel).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.278.23"
$ws.Range("E2").Value = "  +0.00%  "

$ws.Range("D3").Value = "1.930.87"
$ws.Range("E3").Value = "  +0.01%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "0.7525"
$ws.Range("E5").Value = "  +5.07%  "

$ws.Range("D6").Value = "242.29"
$ws.Range("E6").Value = "  -2.72%  "

$ws.Range("D7").Value = "0.9995"
$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").Value = "27.82"
$ws.Range("E8").Value = "  +0.40%  "

$ws.Range("D9").Value = "0.3178"
$ws.Range("E9").Value = "  -0.71%  "

$ws.Range("D10").Value = "0.07109"
$ws.Range("E10").Value = "  +0.07%  "

$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").Value = "0.7797"
$ws.Range("E11").Value = "  -1.56%  "

$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "0.08046"
$ws.Range("E12").Value = "  +0.74%  "

$ws.Range("D13").Value = "1.921.74"
$ws.Range("E13").Value = "  -0.34%  "

$ws.Range("D14").Value = "5.392"
$ws.Range("E14").Value = "  -0.02%  "

$ws.Range("D15").Value = "93.05"
$ws.Range("E15").Value = "  -1.86%  "

$ws.Range("D16").Value = "14.56"
$ws.Range("E16").Value = "  -0.80%  "

$ws.Range("D17").Value = "30.269.08"
$ws.Range("E17").Value = "  +0.04%  "

$ws.Range("D18").Value = "6.010"
$ws.Range("E18").Value = "  +4.04%  "

$ws.Range("D19").Value = "251.96"
$ws.Range("E19").Value = "  -1.65%  "

$ws.Range("D20").Value = "0.000007947"
$ws.Range("E20").Value = "  -1.11%  "

$ws.Range("D21").Value = "2.175.17"
$ws.Range("E21").Value = "  -0.12%  "

$ws.Range("D22").Value = "0.9994"
$ws.Range("E22").Value = "  +0.00%  "

$ws.Range("D23").Value = "0.9990"
$ws.Range("E23").Value = "  -0.05%  "

$ws.Range("D24").Value = "6.689"
$ws.Range("E24").Value = "  -1.94%  "

$ws.Range("D25").Value = "9.546"
$ws.Range("E25").Value = "  +0.06%  "

$ws.Range("D26").Value = "164.94"
$ws.Range("E26").Value = "  -0.06%  "

$ws.Range("D27").Value = "19.11"
$ws.Range("E27").Value = "  +0.02%  "

$ws.Range("D28").Value = "0.1301"
$ws.Range("E28").Value = "  +2.51%  "

$ws.Range("D29").Value = "2.187"
$ws.Range("E29").Value = "  -3.49%  "

$ws.Range("E30").Value = "  +0.86%  "

$ws.Range("D31").Value = "1.546"
$ws.Range("E31").Value = "  +1.18%  "

$ws.Range("D32").Value = "4.413"
$ws.Range("E32").Value = "  +0.43%  "

$ws.Range("D33").Value = "4.145"
$ws.Range("E33").Value = "  +0.30%  "

$ws.Range("D34").Value = "0.05222"

$ws.Range("D35").Value = "1.316"
$ws.Range("E35").Value = "  +3.56%  "

$ws.Range("D36").Value = "0.7586"
$ws.Range("E36").Value = "  +1.78%  "

$ws.Range("D37").Value = "2.782"
$ws.Range("E37").Value = "  +0.34%  "

$ws.Range("D38").Value = "0.01955"
$ws.Range("E38").Value = "  -0.26%  "

$ws.Range("D39").Value = "2.795"
$ws.Range("E39").Value = "  -0.05%  "

$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "6.516"
$ws.Range("E40").Value = "  +2.47%  "

$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "78.35"
$ws.Range("E41").Value = "  -0.51%  "

$ws.Range("D42").Value = "0.4541"
$ws.Range("E42").Value = "  +0.67%  "

$ws.Range("D43").Value = "1.980"
$ws.Range("E43").Value = "  -0.44%  "

$ws.Range("D44").Value = "0.8400"
$ws.Range("E44").Value = "  -0.91%  "

$ws.Range("D45").Value = "0.9992"
$ws.Range("E45").Value = "  -0.03%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "9.994"
$ws.Range("E46").Value = "  +2.60%  "

$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").Value = "7.673"
$ws.Range("E47").Value = "  +3.30%  "

$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "101.59"
$ws.Range("E48").Value = "  +1.11%  "

$ws.Range("D49").Value = "38.00"
$ws.Range("E49").Value = "  +3.67%  "

$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "0.1224"
$ws.Range("E50").Value = "  +7.22%  "

$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "960.27"
$ws.Range("E51").Value = "  +1.26%  "

# Restore the original (default) cell style now that the text values are safely stored,
# so the saved workbook does not pick up a new/extra style index on these cells.
$ws.Range("D2:D51").Style = "Normal"
